$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# I1/J1 header cells, then set their labels.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row-by-row data for the new I0 / IF columns (row, I-value, J-value).
$rowData = @(
    @(2,4,5),
    @(3,8,9),
    @(4,5,5),
    @(5,4,4),
    @(6,6,7),
    @(7,11,11),
    @(8,8,8),
    @(9,11,11),
    @(10,7,7),
    @(11,7,8),
    @(12,6,6),
    @(13,6,6),
    @(14,7,7),
    @(15,8,8),
    @(16,6,6),
    @(17,9,9),
    @(18,9,9),
    @(19,8,8),
    @(20,8,8),
    @(21,6,7),
    @(22,7,7),
    @(23,6,6),
    @(24,7,7),
    @(25,8,8),
    @(26,6,6),
    @(27,6,7),
    @(28,7,7),
    @(29,5,5),
    @(30,6,6),
    @(31,7,7),
    @(32,7,7),
    @(33,6,6),
    @(34,6,6),
    @(35,8,8),
    @(36,6,7),
    @(37,8,8),
    @(38,6,6),
    @(39,7,7),
    @(40,10,10),
    @(41,9,9),
    @(42,5,5),
    @(43,8,8),
    @(44,3,4),
    @(45,8,8),
    @(46,7,7),
    @(47,7,7),
    @(48,8,8),
    @(49,8,8),
    @(50,5,5),
    @(51,7,7),
    @(52,7,7),
    @(53,6,6),
    @(54,9,9),
    @(55,7,7),
    @(56,8,8),
    @(57,7,7),
    @(58,7,7),
    @(59,5,5),
    @(60,7,7),
    @(61,7,7),
    @(62,6,6),
    @(63,8,8),
    @(64,7,7),
    @(65,8,8),
    @(66,8,8),
    @(67,7,8),
    @(68,9,9),
    @(69,10,10),
    @(70,8,8),
    @(71,7,7),
    @(72,6,6),
    @(73,6,6),
    @(74,4,4),
    @(75,6,6),
    @(76,5,5),
    @(77,8,8),
    @(78,6,6),
    @(79,3,3)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
